# Update MS account mappings in Mappings.xlsx.
$wb = $excel.ActiveWorkbook
$ms = $wb.Worksheets.Item("MS")

# Updated account-number mappings (column A) on the MS sheet.
$ms.Range("A3").Value  = "Cardor Alef SpA - 4165"
$ms.Range("A4").Value  = "Cuenta Personal - Jaime - 4120"
$ms.Range("A5").Value  = "NNW Capital SpA - 4105"
$ms.Range("A6").Value  = "NNW II Capital SpA - 4103"
$ms.Range("A7").Value  = "FNW Capital SpA - 4157"
$ms.Range("A10").Value = "Alanseb LP - 4156"
$ms.Range("A11").Value = "Cuenta Personal - Felipe - 4122"
$ms.Range("A12").Value = "Cuenta Personal - Irene - 4162"
$ms.Range("A13").Value = "Cuenta Personal - Jaime y Felipe - 4114"
$ms.Range("A14").Value = "Cuenta Personal - Jaime y Natalia - 4118"
$ms.Range("A15").Value = "Cuenta Personal - Jaime y Nicolas - 4121"
$ms.Range("A16").Value = "Cuenta Personal - Natalia - 4108"
$ms.Range("A17").Value = "Cuenta Personal - Nicolas - 4110"
$ms.Range("A20").Value = "NNW Ventures LLC - 4107"
$ms.Range("A21").Value = "Cuenta Personal - Jacques - 4155"

# Move the active tab / selection from JPM to MS, matching the saved view state.
$ms.Activate()
$ms.Range("H12").Select()
